$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.Style = "Normal"
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "36.467.63"
Set-TextValue "E2" "  +0.13%  "

# Row 3
Set-TextValue "D3" "1.949.20"
Set-TextValue "E3" "  -0.19%  "

# Row 4
Set-TextValue "E4" "  -0.04%  "

# Row 5
Set-TextValue "D5" "243.13"
Set-TextValue "E5" "  -0.39%  "

# Row 6
Set-TextValue "D6" "0.612"
Set-TextValue "E6" "  -0.71%  "

# Row 7
Set-TextValue "B7" "Solana"
Set-TextValue "C7" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D7" "57.92"
Set-TextValue "E7" "  +0.56%  "

# Row 8
Set-TextValue "B8" "USDC"
Set-TextValue "C8" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  -0.10%  "

# Row 9
Set-TextValue "D9" "0.375"
Set-TextValue "E9" "  +2.65%  "

# Row 10
Set-TextValue "D10" "0.0788"
Set-TextValue "E10" "  -7.44%  "

# Row 11
Set-TextValue "E11" "  -1.48%  "

# Row 12
Set-TextValue "D12" "2.234.90"
Set-TextValue "E12" "  -0.23%  "

# Row 13
Set-TextValue "D13" "0.825"
Set-TextValue "E13" "  +0.42%  "

# Row 14
Set-TextValue "D14" "13.71"
Set-TextValue "E14" "  +1.23%  "

# Row 15
Set-TextValue "D15" "21.30"
Set-TextValue "E15" "  -0.86%  "

# Row 16
Set-TextValue "D16" "5.27"
Set-TextValue "E16" "  +1.09%  "

# Row 17
Set-TextValue "D17" "1.948.83"
Set-TextValue "E17" "  -0.08%  "

# Row 18
Set-TextValue "D18" "36.324.84"
Set-TextValue "E18" "  -0.15%  "

# Row 19
Set-TextValue "D19" "69.14"
Set-TextValue "E19" "  -0.94%  "

# Row 20
Set-TextValue "D20" "0.0₃0845"
Set-TextValue "E20" "  -4.41%  "

# Row 21
Set-TextValue "D21" "227.92"
Set-TextValue "E21" "  -0.89%  "

# Row 22
Set-TextValue "D22" "5.01"
Set-TextValue "E22" "  -1.29%  "

# Row 23
Set-TextValue "E23" "  -0.17%  "

# Row 24
Set-TextValue "D24" "2.48"
Set-TextValue "E24" "  +1.99%  "

# Row 25
Set-TextValue "D25" "2.36"
Set-TextValue "E25" "  +2.62%  "

# Row 26
Set-TextValue "E26" "  -2.25%  "

# Row 27
Set-TextValue "B27" "Monero"
Set-TextValue "C27" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D27" "160.42"
Set-TextValue "E27" "  -0.75%  "

# Row 28
Set-TextValue "B28" "Kaspa"
Set-TextValue "C28" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D28" "0.136"
Set-TextValue "E28" "  +1.42%  "

# Row 29
Set-TextValue "E29" "  -0.91%  "

# Row 30
Set-TextValue "D30" "0.120"
Set-TextValue "E30" "  +1.01%  "

# Row 31
Set-TextValue "D31" "1.17"
Set-TextValue "E31" "  +1.20%  "

# Row 32
Set-TextValue "D32" "4.66"
Set-TextValue "E32" "  +0.08%  "

# Row 33
Set-TextValue "E33" "  -6.76%  "

# Row 34
Set-TextValue "D34" "4.31"
Set-TextValue "E34" "  +0.66%  "

# Row 35
Set-TextValue "E35" "  -0.08%  "

# Row 36
Set-TextValue "D36" "3.43"
Set-TextValue "E36" "  +13.08%  "

# Row 37
Set-TextValue "D37" "2.24"
Set-TextValue "E37" "  +3.62%  "

# Row 38
Set-TextValue "E38" "  -1.41%  "

# Row 39
Set-TextValue "D39" "5.25"
Set-TextValue "E39" "  -15.09%  "

# Row 40
Set-TextValue "D40" "0.0969"
Set-TextValue "E40" "  -1.65%  "

# Row 41
Set-TextValue "E41" "  -0.41%  "

# Row 42
Set-TextValue "E42" "  -1.23%  "

# Row 43
Set-TextValue "E43" "  -1.23%  "

# Row 44
Set-TextValue "D44" "15.70"
Set-TextValue "E44" "  -0.26%  "

# Row 45
Set-TextValue "D45" "1.358.58"
Set-TextValue "E45" "  -0.09%  "

# Row 46
Set-TextValue "E46" "  -0.84%  "

# Row 47
Set-TextValue "D47" "87.20"
Set-TextValue "E47" "  -0.67%  "

# Row 48
Set-TextValue "D48" "7.09"
Set-TextValue "E48" "  -0.89%  "

# Row 49
Set-TextValue "D49" "2.83"
Set-TextValue "E49" "  -0.65%  "

# Row 50
Set-TextValue "D50" "2.126.03"
Set-TextValue "E50" "  -0.24%  "

# Row 51
Set-TextValue "D51" "43.56"
Set-TextValue "E51" "  -3.32%  "
